# The commit swaps the presentation's theme ("Integral") for the built-in
# "Office Theme" palette (design/theme gallery change). In the OOXML this
# shows up as the 12 theme colours (ppt/theme/theme1.xml's <a:clrScheme>)
# changing from the Integral palette to the standard Office palette.
#
# Helper: build a COM-style BGR-packed RGB long from R,G,B bytes (same
# encoding PowerPoint's ColorFormat.RGB uses: R + G*256 + B*65536).
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The slide master's Theme drives the deck's design (ppt/theme/theme1.xml).
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# ThemeColorScheme item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1,
# 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink,
# 12 folHlink -- apply the "Office Theme" swatches in that order.
$colors.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$colors.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colors.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$colors.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colors.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colors.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$colors.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colors.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$colors.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$colors.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$colors.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$colors.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72

# Relabel the design from "Integral" to the built-in "Office Theme".
try { $colors.Name = "Office" } catch {}
try { $theme.Name = "Office Theme" } catch {}
